# Commit: Wed, Jun 24, 2020 10:06:13 AM
#
# 1) Slide 5's table switches from the custom "Table_0" style to the
#    built-in "Medium Style 2 - Accent 1" table style.
# 2) The presentation's Design (slideMaster1 -> theme1.xml) is switched
#    from the "Integral" / "Red Violet" color theme back to the plain
#    "Office Theme" / "Office" color theme (the theme that ships with a
#    brand-new PowerPoint deck). This is what happens when a user opens
#    the Design tab and clicks the built-in "Office Theme" swatch in the
#    gallery - PowerPoint rewrites the twelve theme colors used by the
#    slide master/layouts/slides.

$p = $ppt.ActivePresentation

# --- 1. Re-style the table on slide 5 ("B1- TYPES OF FINANCIAL DOCUMENTS") ---
$tableSlide = $p.Slides.Item(5)
$tableShape = $tableSlide.Shapes.Item(2)
$tbl = $tableShape.Table
$tbl.ApplyStyle("{FD8E245D-EDCB-4DA4-890A-52AC85C74F49}", $true)

# --- 2. Re-color the presentation theme to the stock "Office Theme" ---
$themeSlide = $p.Slides.Item(1)
$tcs = $themeSlide.ThemeColorScheme

$tcs.Item(1).RGB  = 0            # dk1      000000
$tcs.Item(2).RGB  = 16777215     # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388      # dk2      44546A
$tcs.Item(4).RGB  = 15132391     # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939     # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501      # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845     # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407        # accent4  FFC000
$tcs.Item(9).RGB  = 12874308     # accent5  4472C4
$tcs.Item(10).RGB = 4697456      # accent6  70AD47
$tcs.Item(11).RGB = 12673797     # hlink    0563C1
$tcs.Item(12).RGB = 7491477      # folHlink 954F72
